$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "保險" (insurance) — add company/name/owner headers + new
# metadata columns E:K (property_category, category, date,
# legislator_name, legislator_id, source_file, index)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("保險")

$ws4.Range("B1").Value = "company"
$ws4.Range("C1").Value = "name"
$ws4.Range("D1").Value = "owner"
$ws4.Range("E1").Value = "property_category"
$ws4.Range("F1").Value = "category"
$ws4.Range("G1").Value = "date"
$ws4.Range("H1").Value = "legislator_name"
$ws4.Range("I1").Value = "legislator_id"
$ws4.Range("J1").Value = "source_file"
$ws4.Range("K1").Value = "index"

$ws4.Range("E2:E4").Value = "insurance"
$ws4.Range("F2:F4").Value = "normal"

$ws4.Range("G2:G4").NumberFormat = "@"
$ws4.Range("G2").Value = "2012-04-30"
$ws4.Range("G3").Value = "2012-04-30"
$ws4.Range("G4").Value = "2012-04-30"
$ws4.Range("G2:G4").Style = "Normal"

$ws4.Range("H2:H4").Value = "蔡錦隆"
$ws4.Range("I2:I4").Value = 1380
$ws4.Range("J2:J4").Value = "tmpf6571"

$ws4.Range("K2").Value = 128
$ws4.Range("K3").Value = 129
$ws4.Range("K4").Value = 130

# ---------------------------------------------------------------------
# Sheet "債務" (debt) — fix row1 header labels (species/debtor/owner/
# total/register_date/register_reason) + new metadata columns H:N
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("債務")

$ws5.Range("B1").Value = "species"
$ws5.Range("C1").Value = "debtor"
$ws5.Range("D1").Value = "owner"
$ws5.Range("E1").Value = "total"
$ws5.Range("F1").Value = "register_date"
$ws5.Range("G1").Value = "register_reason"
$ws5.Range("H1").Value = "property_category"
$ws5.Range("I1").Value = "category"
$ws5.Range("J1").Value = "date"
$ws5.Range("K1").Value = "legislator_name"
$ws5.Range("L1").Value = "legislator_id"
$ws5.Range("M1").Value = "source_file"
$ws5.Range("N1").Value = "index"

$ws5.Range("H2:H3").Value = "debt"
$ws5.Range("I2:I3").Value = "normal"

$ws5.Range("J2:J3").NumberFormat = "@"
$ws5.Range("J2").Value = "2012-04-30"
$ws5.Range("J3").Value = "2012-04-30"
$ws5.Range("J2:J3").Style = "Normal"

$ws5.Range("K2:K3").Value = "蔡錦隆"
$ws5.Range("L2:L3").Value = 1380
$ws5.Range("M2:M3").Value = "tmpf6571"

$ws5.Range("N2").Value = 140
$ws5.Range("N3").Value = 141

# ---------------------------------------------------------------------
# Sheet "事業投資" (investment) — fix row1 header labels (owner/company/
# address/total/register_date/register_reason) + new metadata columns
# H:N
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("事業投資")

$ws6.Range("B1").Value = "owner"
$ws6.Range("C1").Value = "company"
$ws6.Range("D1").Value = "address"
$ws6.Range("E1").Value = "total"
$ws6.Range("F1").Value = "register_date"
$ws6.Range("G1").Value = "register_reason"
$ws6.Range("H1").Value = "property_category"
$ws6.Range("I1").Value = "category"
$ws6.Range("J1").Value = "date"
$ws6.Range("K1").Value = "legislator_name"
$ws6.Range("L1").Value = "legislator_id"
$ws6.Range("M1").Value = "source_file"
$ws6.Range("N1").Value = "index"

$ws6.Range("H2").Value = "investment"
$ws6.Range("I2").Value = "normal"

$ws6.Range("J2").NumberFormat = "@"
$ws6.Range("J2").Value = "2012-04-30"
$ws6.Range("J2").Style = "Normal"

$ws6.Range("K2").Value = "蔡錦隆"
$ws6.Range("L2").Value = 1380
$ws6.Range("M2").Value = "tmpf6571"
$ws6.Range("N2").Value = 146
